$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# New log rows captured by the mmWave sensor
$rows = @(
    @{ Row = 14; Date = "2026-02-01"; Time = "14:41:03"; Hour = "14:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" },
    @{ Row = 15; Date = "2026-02-01"; Time = "14:41:05"; Hour = "14:00"; Location = "Living Room"; Value = "PRESENCE_DETECTED"; Status = "Active" }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $dateCell = $ws.Cells.Item($rowNum, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $r.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($rowNum, 2).Value = $r.Time
    $ws.Cells.Item($rowNum, 3).Value = $r.Hour
    $ws.Cells.Item($rowNum, 4).Value = $r.Location
    $ws.Cells.Item($rowNum, 5).Value = $r.Value
    $ws.Cells.Item($rowNum, 6).Value = $r.Status
}
